# Updates cryptos list figures (Price / Volume(1h) columns, plus a few
# reordered coin rows) to match the latest scrape.
# Note: some Price values look like plain numbers (e.g. '72.00', '9.50').
# Excel auto-converts such text to a number (dropping trailing zeros), so
# those are entered with a leading apostrophe to force them to stay text,
# exactly like the source data (e.g. '''72.00' -> literal text "'72.00").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '36.149.90'
$ws.Range("E2").Value = '  -1.40%  '

# Row 3
$ws.Range("D3").Value = '2.010.47'
$ws.Range("E3").Value = '  -1.98%  '

# Row 4
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").Value = '''252.02'
$ws.Range("E5").Value = '  +2.23%  '

# Row 6
$ws.Range("D6").Value = '''0.643'
$ws.Range("E6").Value = '  -3.55%  '

# Row 7
$ws.Range("D7").Value = '''62.52'
$ws.Range("E7").Value = '  +8.83%  '

# Row 8
$ws.Range("E8").Value = '  +0.03%  '

# Row 9
$ws.Range("D9").Value = '''59.07'
$ws.Range("E9").Value = '  -6.67%  '

# Row 10
$ws.Range("E10").Value = '  -0.60%  '

# Row 11
$ws.Range("D11").Value = '''0.0746'
$ws.Range("E11").Value = '  -1.28%  '

# Row 12
$ws.Range("E12").Value = '  -1.74%  '

# Row 13
$ws.Range("E13").Value = '  -2.21%  '

# Row 14
$ws.Range("D14").Value = '''14.84'
$ws.Range("E14").Value = '  +2.41%  '

# Row 15
$ws.Range("D15").Value = '2.305.11'
$ws.Range("E15").Value = '  -1.92%  '

# Row 16
$ws.Range("D16").Value = '''5.42'
$ws.Range("E16").Value = '  -0.71%  '

# Row 17
$ws.Range("D17").Value = '''19.52'
$ws.Range("E17").Value = '  +8.49%  '

# Row 18
$ws.Range("D18").Value = '2.004.71'
$ws.Range("E18").Value = '  -2.58%  '

# Row 19
$ws.Range("D19").Value = '36.071.18'
$ws.Range("E19").Value = '  -1.28%  '

# Row 20
$ws.Range("D20").Value = '''72.00'
$ws.Range("E20").Value = '  +0.02%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0856'
$ws.Range("E21").Value = '  -0.80%  '

# Row 22
$ws.Range("D22").Value = '''5.26'
$ws.Range("E22").Value = '  +0.53%  '

# Row 23
$ws.Range("D23").Value = '''233.53'
$ws.Range("E23").Value = '  -1.66%  '

# Row 24
$ws.Range("D24").Value = '''2.68'
$ws.Range("E24").Value = '  +17.68%  '

# Row 25
$ws.Range("E25").Value = '  +0.05%  '

# Row 26
$ws.Range("E26").Value = '  -2.90%  '

# Row 27
$ws.Range("D27").Value = '''9.50'
$ws.Range("E27").Value = '  +1.25%  '

# Row 28
$ws.Range("D28").Value = '''163.92'
$ws.Range("E28").Value = '  -0.59%  '

# Row 29
$ws.Range("E29").Value = '  -2.61%  '

# Row 30
$ws.Range("E30").Value = '  -1.37%  '

# Row 31
$ws.Range("D31").Value = '''5.11'
$ws.Range("E31").Value = '  +1.47%  '

# Row 32
$ws.Range("E32").Value = '  -1.61%  '

# Row 33
$ws.Range("E33").Value = '  +24.37%  '

# Row 34
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '''4.51'
$ws.Range("E34").Value = '  +1.13%  '

# Row 35
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").Value = '''2.50'
$ws.Range("E35").Value = '  +12.59%  '

# Row 36
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '''0.0603'
$ws.Range("E36").Value = '  +0.09%  '

# Row 37
$ws.Range("E37").Value = '  -0.11%  '

# Row 38
$ws.Range("E38").Value = '  -1.21%  '

# Row 39
$ws.Range("E39").Value = '  +14.50%  '

# Row 40
$ws.Range("E40").Value = '  +13.76%  '

# Row 41
$ws.Range("D41").Value = '''1.22'
$ws.Range("E41").Value = '  -1.38%  '

# Row 42
$ws.Range("E42").Value = '  +0.02%  '

# Row 43
$ws.Range("E43").Value = '  -1.15%  '

# Row 44
$ws.Range("E44").Value = '  +1.08%  '

# Row 45
$ws.Range("D45").Value = '''16.61'
$ws.Range("E45").Value = '  +3.41%  '

# Row 46
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '''7.80'
$ws.Range("E46").Value = '  +5.15%  '

# Row 47
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '''93.66'
$ws.Range("E47").Value = '  -0.82%  '

# Row 48
$ws.Range("D48").Value = '1.420.81'
$ws.Range("E48").Value = '  +2.75%  '

# Row 49
$ws.Range("D49").Value = '''2.48'
$ws.Range("E49").Value = '  +8.65%  '

# Row 50
$ws.Range("E50").Value = '  -1.23%  '

# Row 51
$ws.Range("D51").Value = '''47.50'
$ws.Range("E51").Value = '  +3.03%  '
